$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-7
# from 2023-10-05 (45204) to 2023-10-08 (45207)
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45207
}
